# Apply updated "dSF" (column F) values.
# These come from a repull/recalculation of data where the end-of-season
# snapshot value (dSF) differs from the season-start value (dS0) captured
# in column E for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -2
    5  = -3
    9  = -3
    10 = 2
    14 = 1
    20 = -2
    32 = -1
    39 = 4
    41 = -3
    43 = 3
    50 = 1
    51 = -4
    55 = -5
    61 = -2
    62 = 4
    65 = -2
    66 = -9
    70 = -1
    71 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
